# Finish laravel tutorial_05 review row (row 9) on the "Review" sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review")

# No. / date
$ws.Cells.Item(9, 1).Value = 8
$ws.Cells.Item(9, 2).Value = 44572

# Status / part / error type (reuse existing values, same as row 8)
$ws.Cells.Item(9, 3).Value = "Open"
$ws.Cells.Item(9, 4).Value = "Others"
$ws.Cells.Item(9, 5).Value = "Careless"

# Location / review content
$ws.Cells.Item(9, 6).Value = "assignment_05 and all "
$ws.Cells.Item(9, 9).Value = "1)use laravel with() 2)write function doc 3)write db code in DAO 4)use form request 5)use DB transaction 6)add else condition"

# Reviewer / response
$ws.Cells.Item(9, 15).Value = "PyaePyaeHan"
$ws.Cells.Item(9, 16).Value = "complete"

# Response date - copy the date number format from the row above (V8) so the
# new date (V9) renders the same as the rest of the column, then set value.
$ws.Cells.Item(8, 22).Copy()
$ws.Cells.Item(9, 22).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(9, 22).Value = 44572

# Response person
$ws.Cells.Item(9, 23).Value = "WaiLinOo"

# Move the active selection to X9 (matches cursor position after filling the row)
[void]$ws.Range("X9").Select()
